$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column BJ (Visual_exp)
$ws.Range("BJ1").Value = "Visual_exp"

# Values for BJ2:BJ78 (Visual_exp data per row)
$visualExp = @{
    2 = 105
    3 = 0
    4 = 155
    5 = 243
    6 = 155
    7 = 1108
    8 = 0
    9 = 0
    10 = 245
    11 = 0
    12 = 0
    13 = 140
    14 = 0
    15 = 0
    16 = 0
    17 = 110
    18 = 294
    19 = 0
    20 = 123
    21 = 531
    22 = 500
    23 = 485
    24 = 237
    25 = 0
    26 = 0
    27 = 289
    28 = 205
    29 = 253
    30 = 193
    31 = 359
    32 = 407
    33 = 0
    34 = 743
    35 = 350
    36 = 371
    37 = 0
    38 = 68
    39 = 172
    40 = 447
    41 = 0
    42 = 136
    43 = 0
    44 = 549
    45 = 132
    46 = 106
    47 = 191
    48 = 211
    49 = 525
    50 = 89
    51 = 571
    52 = 419
    53 = 64
    54 = 138
    55 = 690
    56 = 281
    57 = 0
    58 = 177
    59 = 433
    60 = 311
    61 = 511
    62 = 162
    63 = 572
    64 = 0
    65 = 201
    66 = 656
    67 = 221
    68 = 0
    69 = 586
    70 = 532
    71 = 254
    72 = 0
    73 = 884
    74 = 256
    75 = 255
    76 = 521
    77 = 1806
    78 = 83
}

foreach ($row in $visualExp.Keys) {
    $ws.Cells.Item($row, 62).Value = $visualExp[$row]
}
